$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.540.33"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "2.412.78"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'507.76"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").Value = "'132.75"
$ws.Range("E6").Value = "  +1.64%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").Value = "2.451.36"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'0.0984"
$ws.Range("E10").Value = "  +0.31%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "'4.64"
$ws.Range("E13").Value = "  -5.87%  "
$ws.Range("D14").Value = "2.846.74"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "57.411.01"
$ws.Range("D16").Value = "'21.95"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "2.431.71"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'10.33"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'314.42"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").Value = "'6.35"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "'5.68"
$ws.Range("E24").Value = "  -2.49%  "
$ws.Range("D25").Value = "'65.36"
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.548.92"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'0.993"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'0.383"
$ws.Range("E28").Value = "  -4.64%  "
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("D30").Value = "'7.63"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").Value = "'174.28"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "0.0₃0738"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'1.69"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "'6.22"
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'0.991"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").Value = "'18.00"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "'1.24"
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("D40").Value = "'3.89"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "'0.821"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").Value = "'36.50"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "'1.47"
$ws.Range("E43").Value = "  +1.42%  "
$ws.Range("D44").Value = "'134.34"
$ws.Range("E44").Value = "  +9.45%  "
$ws.Range("D45").Value = "'3.42"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'4.91"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").Value = "'259.48"
$ws.Range("E47").Value = "  -1.54%  "
$ws.Range("D48").Value = "'0.572"
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("D49").Value = "'0.0919"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").Value = "'0.0496"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  +1.96%  "
